# Updated cryptos list with latest price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> updated Price (column D) / Volume(1h) (column E) text.
$updates = @{
    2 = @{D='25.991.61'; E='  +0.64%  '}
    3 = @{D='1.745.93'; E='  -0.48%  '}
    4 = @{E='  +0.09%  '}
    5 = @{D='232.90'; E='  -1.74%  '}
    6 = @{E='  +0.10%  '}
    7 = @{D='0.5187'; E='  +1.75%  '}
    8 = @{D='0.2818'; E='  +4.35%  '}
    9 = @{D='39.58'; E='  -3.13%  '}
    10 = @{D='0.06127'; E='  -1.28%  '}
    11 = @{D='1.753.38'; E='  +0.01%  '}
    12 = @{D='0.07006'; E='  +1.20%  '}
    13 = @{D='15.38'; E='  -1.28%  '}
    14 = @{D='0.6411'; E='  +1.83%  '}
    15 = @{D='4.518'; E='  +0.43%  '}
    16 = @{D='77.40'; E='  -1.32%  '}
    17 = @{E='  +0.11%  '}
    19 = @{D='25.996.40'; E='  +0.57%  '}
    21 = @{D='0.000006603'; E='  -1.86%  '}
    22 = @{D='1.979.80'; E='  +0.50%  '}
    23 = @{E='  +1.66%  '}
    24 = @{D='8.639'; E='  +4.40%  '}
    25 = @{D='5.135'; E='  -0.99%  '}
    26 = @{D='139.17'; E='  +1.85%  '}
    27 = @{D='1.503'; E='  +2.84%  '}
    28 = @{D='15.09'; E='  -1.69%  '}
    29 = @{D='1.820'; E='  +1.67%  '}
    30 = @{D='102.49'; E='  -0.30%  '}
    31 = @{E='  -0.04%  '}
    32 = @{D='3.659'; E='  -1.94%  '}
    33 = @{E='  +0.15%  '}
    34 = @{D='0.04469'; E='  +1.19%  '}
    35 = @{D='2.610'; E='  -1.47%  '}
    36 = @{D='0.9893'; E='  -1.70%  '}
    37 = @{D='0.6152'; E='  +1.51%  '}
    38 = @{D='2.672'; E='  -0.86%  '}
    39 = @{D='0.01591'}
    40 = @{D='1.928'; E='  -2.19%  '}
    41 = @{D='1.000'; E='  -0.01%  '}
    42 = @{D='100.25'; E='  -1.99%  '}
    43 = @{D='0.3856'; E='  -0.27%  '}
    44 = @{D='0.7452'; E='  -0.80%  '}
    45 = @{D='5.065'; E='  +3.27%  '}
    46 = @{D='0.05462'; E='  -0.85%  '}
    47 = @{D='6.325'; E='  +5.79%  '}
    48 = @{D='0.1126'; E='  +2.58%  '}
    49 = @{D='53.03'; E='  +0.20%  '}
    50 = @{D='30.04'; E='  -0.72%  '}
    51 = @{D='7.623'; E='  +2.12%  '}
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    if ($vals.ContainsKey('D')) {
        # Some "prices" look like numbers (e.g. 232.90, 0.5187) but must stay
        # literal text (matching the sheet's existing inline-string cells,
        # multi-dot thousand groupings, trailing zeros, etc.). Force text
        # format while assigning, then drop the format again so the cell's
        # style stays at the sheet default, same as before the edit.
        $ws.Range("D$row").NumberFormat = "@"
        $ws.Range("D$row").Value = $vals['D']
        $ws.Range("D$row").ClearFormats()
    }

    if ($vals.ContainsKey('E')) {
        # "  +0.64%  " etc. already has padding spaces so Excel keeps it text.
        $ws.Range("E$row").Value = $vals['E']
    }
}
